$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (row 1) holds two parallel blocks of ten columns each,
# describing the same set of fields once for the "old" input file and
# once for the "new" input file, plus a "diff" column in between.
# Rename the "_old" / "_new" suffixes to the actual format-version
# identifiers that the two input files represent.
$fieldNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $fieldNames.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $fieldNames[$i] + "_FV2310"
    $ws.Range($newCols[$i] + "1").Value = $fieldNames[$i] + "_FV2404"
}

# Turn the used range into an actual Excel table (ListObject) with the
# freshly renamed headers, auto filter included.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U73"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row so it stays visible while scrolling.
$null = $ws.Range("A2").Select()
$null = ($excel.ActiveWindow.FreezePanes = $true)
